$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outputs")

# Insert a new row at position 2 (before "general"), shifting existing rows down
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with the "raw" entry
$ws.Range("A2").Value = "raw"
$ws.Range("B2").Value = 1

# Rename the "perfectionism/certainty" label (now in row 4) to "perfectionism_certainty"
$ws.Range("A4").Value = "perfectionism_certainty"
